# Version with a working long press button
# Adds 12 new rows (17-28) to the "Translation" sheet describing the
# text ids used for the new long-press button behaviour (Delay / Duration
# in ms, aligned Left/Right), mirroring the rows already present in the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$rows = @(
    @("SingleUseId17", "Default", "Right", "LTR", "Delay",    "New Text"),
    @("SingleUseId18", "Default", "Right", "LTR", "Duration", "New Text"),
    @("SingleUseId19", "Default", "Left",  "LTR", "ms",       "New Text"),
    @("SingleUseId20", "Default", "Left",  "LTR", "ms",       "New Text"),
    @("SingleUseId29", "Default", "Right", "LTR", "Delay",    "New Text"),
    @("SingleUseId30", "Default", "Right", "LTR", "Duration", "New Text"),
    @("SingleUseId31", "Default", "Left",  "LTR", "ms",       "New Text"),
    @("SingleUseId32", "Default", "Left",  "LTR", "ms",       "New Text"),
    @("SingleUseId33", "Default", "Right", "LTR", "Delay",    "New Text"),
    @("SingleUseId34", "Default", "Right", "LTR", "Duration", "New Text"),
    @("SingleUseId35", "Default", "Left",  "LTR", "ms",       "New Text"),
    @("SingleUseId36", "Default", "Left",  "LTR", "ms",       "New Text")
)

$startRow = 17
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $rowValues[0]  # B: TEXT ID
    $ws.Cells.Item($r, 3).Value = $rowValues[1]  # C: TYPOGRAPHY NAME
    $ws.Cells.Item($r, 4).Value = $rowValues[2]  # D: ALIGNMENT
    $ws.Cells.Item($r, 5).Value = $rowValues[3]  # E: DIRECTION
    $ws.Cells.Item($r, 6).Value = $rowValues[4]  # F: GB
    $ws.Cells.Item($r, 7).Value = $rowValues[5]  # G: SI
}
